$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to unify the DataNode/DataTable/Entity concept (per commit message)
$ws.Name = "DataNode"

# Minor column-width touch-ups (as recorded when the workbook was resaved)
$ws.Columns.Item(1).ColumnWidth = 11.91
$ws.Columns.Item(5).ColumnWidth = 13.16

# Selection moved before the file was saved
$ws.Range("E23").Select()
